$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds plain-text numeric-looking strings (e.g. "1.003",
# "30.866.97" with thousands separators as literal dots). Force the column
# to Text format first so COM does not coerce these into floating point
# numbers (which would lose the exact original formatting/precision).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.866.97"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "2.114.01"
$ws.Range("E3").Value = "  +9.78%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "335.50"
$ws.Range("E5").Value = "  +4.94%  "
$ws.Range("D7").Value = "0.5304"
$ws.Range("E7").Value = "  +4.37%  "
$ws.Range("D8").Value = "0.4360"
$ws.Range("E8").Value = "  +7.96%  "
$ws.Range("D9").Value = "0.08994"
$ws.Range("E9").Value = "  +7.94%  "
$ws.Range("D10").Value = "45.86"
$ws.Range("E10").Value = "  +8.46%  "
$ws.Range("D11").Value = "1.177"
$ws.Range("E11").Value = "  +5.32%  "
$ws.Range("D12").Value = "25.01"
$ws.Range("E12").Value = "  +3.35%  "
$ws.Range("D13").Value = "2.110.70"
$ws.Range("E13").Value = "  +9.59%  "
$ws.Range("D14").Value = "6.757"
$ws.Range("E14").Value = "  +5.17%  "
$ws.Range("D15").Value = "7.794"
$ws.Range("E15").Value = "  +7.32%  "
$ws.Range("D16").Value = "97.43"
$ws.Range("E16").Value = "  +5.09%  "
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "0.00001134"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").Value = "0.06664"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "19.08"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "6.360"
$ws.Range("E22").Value = "  +6.71%  "
$ws.Range("D23").Value = "30.941.61"
$ws.Range("E23").Value = "  +2.72%  "
$ws.Range("D24").Value = "12.14"
$ws.Range("E24").Value = "  +6.90%  "
$ws.Range("D25").Value = "2.361.19"
$ws.Range("E25").Value = "  +10.36%  "
$ws.Range("D26").Value = "2.271"
$ws.Range("D27").Value = "22.75"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").Value = "2.564"
$ws.Range("E28").Value = "  +12.43%  "
$ws.Range("D29").Value = "162.77"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "133.24"
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").Value = "1.168"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").Value = "6.230"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("D34").Value = "4.014"
$ws.Range("E34").Value = "  +6.06%  "
$ws.Range("D35").Value = "1.526"
$ws.Range("E35").Value = "  +22.74%  "
$ws.Range("D36").Value = "0.02614"
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("D37").Value = "5.537"
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("D38").Value = "12.86"
$ws.Range("E38").Value = "  +11.02%  "
$ws.Range("D39").Value = "9.533"
$ws.Range("E39").Value = "  +9.99%  "
$ws.Range("D40").Value = "0.06726"
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("D41").Value = "0.2270"
$ws.Range("E41").Value = "  +5.62%  "
$ws.Range("D42").Value = "0.6847"
$ws.Range("E42").Value = "  +5.61%  "
$ws.Range("D43").Value = "1.250"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.6455"
$ws.Range("E44").Value = "  +6.56%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "14.09"
$ws.Range("E46").Value = "  +5.93%  "
$ws.Range("D47").Value = "2.237"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").Value = "3.691"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "1.277"
$ws.Range("D50").Value = "82.51"
$ws.Range("E50").Value = "  +5.72%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "119.80"
$ws.Range("E51").Value = "  -2.19%  "
